$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 749, shifting existing rows 749:825 down to 750:826
$ws.Rows.Item(749).Insert()

# Populate the newly inserted row 749 with the new data record
$ws.Cells.Item(749, 1).Value = 6
$ws.Cells.Item(749, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(749, 3).Value = "Metropolitana"
$ws.Cells.Item(749, 4).Value = 45194
$ws.Cells.Item(749, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(749, 5).Value = 13
$ws.Cells.Item(749, 6).Value = 100112012
$ws.Cells.Item(749, 7).Value = "Espinaca"
$ws.Cells.Item(749, 8).Value = "Sin especificar"
$ws.Cells.Item(749, 9).Value = "Primera"
$ws.Cells.Item(749, 10).Value = 280
$ws.Cells.Item(749, 11).Value = 7000
$ws.Cells.Item(749, 12).Value = 8000
$ws.Cells.Item(749, 13).Value = 7536
$ws.Cells.Item(749, 14).Value = "$/cuna 10 kilos"
$ws.Cells.Item(749, 15).Value = "Región Metropolitana"
$ws.Cells.Item(749, 16).Value = 754
$ws.Cells.Item(749, 17).Value = 10
$ws.Cells.Item(749, 18).Value = "Hortaliza"
